# TC13_Bento_Filter_Chemo-TAC.xlsx — "Fixed Bento 80 Test scripts"
#
# The three Cypher queries on the "startup" sheet (Samples tab query in B2,
# Files tab query in B3, and the file-list query in B4) each get an
# " order By ... ASC LIMIT 100" clause appended/applied, the affected rows
# grow by one wrapped line, and the previously-selected cell moves from C4
# to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: SamplesTab query (row 2) ---------------------------------------
$samplesQueryOld = $ws.Range("B2").Value()
$samplesQueryNew = $samplesQueryOld + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Range("B2").Value = $samplesQueryNew

# --- B3: FilesTab query (row 3) -----------------------------------------
$filesQueryOld = $ws.Range("B3").Value()
$filesQueryNew = $filesQueryOld + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value = $filesQueryNew

# --- B4: file listing query (row 4) --------------------------------------
$lastQueryOld = $ws.Range("B4").Value()
$lastQueryNew = $lastQueryOld -replace "`n    order by f\.file_name$", "`n     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value = $lastQueryNew

# --- Row heights grow by one wrapped line for rows 2 & 3 (row 4 is already
#     clamped at Excel's 409.6pt maximum, so it stays put). ----------------
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# --- Selection moves from C4 to B4, with the view scrolled down a bit ----
$ws.Range("B4").Select()
